$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted numbers (e.g. thousands-separated prices such as
# "41.528.25"). For any new value that Excel would otherwise auto-detect as a real
# number, force the cell format to Text first so the stored value stays a string,
# matching the source data (inline/shared strings) exactly.

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.528.25"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").Value = "2.484.87"
$ws.Range("E3").Value = "  +0.78%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "313.40"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("D6").Value = "93.23"
$ws.Range("E6").Value = "  -1.11%  "

$ws.Range("E7").Value = "  -0.79%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("E9").Value = "  -0.79%  "

$ws.Range("E10").Value = "  -3.08%  "

$ws.Range("E11").Value = "  +0.55%  "

$ws.Range("E12").Value = "  +1.89%  "

$ws.Range("D13").Value = "2.868.74"
$ws.Range("E13").Value = "  +0.98%  "

$ws.Range("D14").Value = "6.84"
$ws.Range("E14").Value = "  -1.62%  "

$ws.Range("D15").Value = "15.59"
$ws.Range("E15").Value = "  +6.68%  "

$ws.Range("D16").Value = "2.495.30"
$ws.Range("E16").Value = "  +1.10%  "

$ws.Range("D17").Value = "0.754"
$ws.Range("E17").Value = "  -4.25%  "

$ws.Range("D18").Value = "41.564.85"
$ws.Range("E18").Value = "  +0.78%  "

$ws.Range("D19").Value = "6.35"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("E20").Value = "  +1.18%  "

$ws.Range("E21").Value = "  +4.44%  "

$ws.Range("D22").Value = "11.19"
$ws.Range("E22").Value = "  -2.56%  "

$ws.Range("D23").Value = "236.22"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "2.71"
$ws.Range("E24").Value = "  -2.50%  "

$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("E26").Value = "  -1.13%  "

$ws.Range("D27").Value = "24.82"
$ws.Range("E27").Value = "  +1.68%  "

$ws.Range("D29").Value = "9.66"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("D30").Value = "36.33"
$ws.Range("E30").Value = "  +0.66%  "

$ws.Range("D31").Value = "157.07"
$ws.Range("E31").Value = "  +2.72%  "

$ws.Range("D32").Value = "5.43"
$ws.Range("E32").Value = "  -2.13%  "

$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "2.57"
$ws.Range("E33").Value = "  -1.15%  "

$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "18.16"
$ws.Range("E34").Value = "  +6.27%  "

$ws.Range("D35").Value = "0.0756"
$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("D36").Value = "2.46"
$ws.Range("E36").Value = "  -5.28%  "

$ws.Range("E37").Value = "  -1.92%  "

$ws.Range("D38").Value = "0.106"
$ws.Range("E38").Value = "  +3.00%  "

$ws.Range("E39").Value = "  -3.11%  "

$ws.Range("E40").Value = "  -0.18%  "

$ws.Range("D41").Value = "4.13"
$ws.Range("E41").Value = "  -2.87%  "

$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("D43").Value = "19.83"
$ws.Range("E43").Value = "  -5.84%  "

$ws.Range("D44").Value = "1.963.88"
$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").Value = "0.0285"
$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("D46").Value = "2.97"
$ws.Range("E46").Value = "  -3.18%  "

$ws.Range("D47").Value = "8.85"
$ws.Range("E47").Value = "  +2.14%  "

$ws.Range("D48").Value = "2.727.02"

$ws.Range("D49").Value = "96.60"
$ws.Range("E49").Value = "  -0.75%  "

$ws.Range("D50").Value = "67.49"
$ws.Range("E50").Value = "  -3.48%  "

$ws.Range("D51").Value = "73.47"
$ws.Range("E51").Value = "  -3.50%  "
